$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (standalone "519033 - Carlos Yujiro Shigue" values in B/C,
# with no label in A) is removed entirely; everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# After the shift, several of the "value" cells (columns B/C) now hold content
# that effectively comes from one slot further down the original sequence,
# while the old long-form texts for those slots are dropped. Apply the
# resulting content directly to match the target state.

# Row 10 "Objetivos:" now shows the docente's info instead of the objectives text.
$ws.Range("B10:C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 (was old row 14) "Programa resumido:" now shows "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (was old row 16) "Programa:" now shows the activation date.
$ws.Range("B15:C15").Value = "01/01/2012"

# Row 18 (was old row 19) "Método:" now shows the docente's info.
$ws.Range("B18:C18").Value = "519033 - Carlos Yujiro Shigue"

# Row 19 (was old row 20) "Critério:" now shows the teaching method text.
$ws.Range("B19:C19").Value = "Aulas expositivas, seminários e exercícios comentados."

# Row 20 (was old row 21) "Norma de recuperação:" now shows the grading criterion text.
$ws.Range("B20:C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# Row 21 (was old row 22) "Bibliografia:" now shows the recovery-exam norm text.
$ws.Range("B21:C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
